# Add two new "explanation" strings to the "2_Vocab" worksheet (the traffic-flow
# matching exercise) in column D, next to the existing A:C matching columns.
#
#   D2 -> "This would be beyond the control of the model-maker, so would be a
#          defined parameter."
#   D6 -> "This describes the state of the system at a given moment."
#
# Both new cells reuse the same wrap-text style ("s=1") already used by every
# other cell on the sheet, so no explicit style assignment is required - just
# writing .Value into a previously-empty cell on this sheet inherits it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2_Vocab")

$ws.Range("D2").Value = "This would be beyond the control of the model-maker, so would be a defined parameter."
$ws.Range("D6").Value = "This describes the state of the system at a given moment."

# Row 2 grows taller to accommodate the newly-wrapped text in D2 (Excel
# recalculates the row's natural height once content makes it wrap to three
# lines).
$ws.Rows.Item(2).RowHeight = 45

# Leave the selection on the last-edited cell, matching the saved cursor
# position.
$ws.Range("D6").Select()
